$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Date: update timestamp (kept as text, not a date value)
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive: false -> true
# Assign via a leading apostrophe so Excel stores literal text "true"
# rather than auto-converting it to a Boolean, then restore the original
# (non quote-prefixed) cell formatting from the neighboring cell.
$caseSensitiveCell = $ws.Range("B17")
$caseSensitiveCell.Value = "'true"
$ws.Range("B16").Copy() | Out-Null
$caseSensitiveCell.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
